$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (matches source data which uses
# ambiguous multi-dot / leading-zero numeric strings) so Excel does not
# silently coerce them into numbers and drop formatting (e.g. trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.374.10'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '1.708.28'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('D4').Value = '0.9954'
$ws.Range('D5').Value = '238.80'
$ws.Range('E5').Value = '  -3.26%  '
$ws.Range('D6').Value = '0.9960'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = '0.4876'
$ws.Range('E7').Value = '  -0.92%  '
$ws.Range('D8').Value = '0.2577'
$ws.Range('E8').Value = '  -3.98%  '
$ws.Range('D9').Value = '0.06159'
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('D10').Value = '1.715.59'
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('D11').Value = '0.06947'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').Value = '15.45'
$ws.Range('E12').Value = '  -1.71%  '
$ws.Range('D13').Value = '0.5956'
$ws.Range('E13').Value = '  -3.10%  '
$ws.Range('D14').Value = '4.443'
$ws.Range('E14').Value = '  -3.16%  '
$ws.Range('D15').Value = '76.35'
$ws.Range('E15').Value = '  -2.29%  '
$ws.Range('D16').Value = '0.9964'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = '26.260.75'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '0.9961'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').Value = '0.000007047'
$ws.Range('E19').Value = '  -3.65%  '
$ws.Range('D20').Value = '11.18'
$ws.Range('E20').Value = '  -3.21%  '
$ws.Range('D21').Value = '1.935.50'
$ws.Range('E21').Value = '  -1.28%  '
$ws.Range('D22').Value = '4.362'
$ws.Range('E22').Value = '  -5.11%  '
$ws.Range('D23').Value = '8.367'
$ws.Range('E23').Value = '  -4.11%  '
$ws.Range('D24').Value = '5.003'
$ws.Range('E24').Value = '  -4.83%  '
$ws.Range('D25').Value = '136.19'
$ws.Range('E25').Value = '  -2.25%  '
$ws.Range('D26').Value = '15.13'
$ws.Range('E26').Value = '  -2.16%  '
$ws.Range('D27').Value = '1.402'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('D28').Value = '1.719'
$ws.Range('E28').Value = '  -2.49%  '
$ws.Range('D29').Value = '105.12'
$ws.Range('E29').Value = '  -2.36%  '
$ws.Range('D30').Value = '3.868'
$ws.Range('E30').Value = '  -4.46%  '
$ws.Range('D31').Value = '0.07924'
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('D32').Value = '3.586'
$ws.Range('E32').Value = '  -3.86%  '
$ws.Range('D33').Value = '0.04438'
$ws.Range('E33').Value = '  -3.98%  '
$ws.Range('D34').Value = '2.601'
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('D35').Value = '0.9904'
$ws.Range('E35').Value = '  -2.77%  '
$ws.Range('D36').Value = '0.6119'
$ws.Range('E36').Value = '  -4.25%  '
$ws.Range('D37').Value = '0.9449'
$ws.Range('E37').Value = '  +4.99%  '
$ws.Range('D38').Value = '1.979'
$ws.Range('E38').Value = '  -3.27%  '
$ws.Range('D39').Value = '2.361'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('D40').Value = '0.9954'
$ws.Range('D41').Value = '0.01474'
$ws.Range('E41').Value = '  -2.00%  '
$ws.Range('D42').Value = '99.33'
$ws.Range('E42').Value = '  -2.59%  '
$ws.Range('D43').Value = '5.369'
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').Value = '0.3782'
$ws.Range('E44').Value = '  -3.68%  '
$ws.Range('D45').Value = '6.780'
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range('D46').Value = '0.1143'
$ws.Range('E46').Value = '  -3.63%  '
$ws.Range('D47').Value = '0.05335'
$ws.Range('E47').Value = '  -1.17%  '
$ws.Range('D48').Value = '30.33'
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').Value = '7.696'
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('D50').Value = '50.85'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('D51').Value = '0.9986'
$ws.Range('E51').Value = '  -0.43%  '
